# edit.ps1
# Applies the "Bai 3: To mau hinh co ban" changes:
#   1. Slide 6 content placeholder gains two descriptive bullet paragraphs
#      about drawing circles.
#   2. Slide 7 title changes from "Ve hinh vong tron" to "To mau hinh tron".
#   3. Slide 7 content placeholder gains two descriptive bullet paragraphs
#      about drawing line segments.
#   4. The cached "datetimeFigureOut" field text on every slide layout and
#      the slide master is bumped from 3/1/2023 to 3/2/2023.

$p = $ppt.ActivePresentation

function Set-TwoParagraphs($shape, $first, $second) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = $first
    $tr.LanguageID = "en-US"
    $tr.InsertAfter("`r" + $second)
    $p2 = $tr.Paragraphs(2, 1)
    $p2.LanguageID = "en-US"
}

# --- Slide 6: "5. To mau hinh tron (Cach 1)" content placeholder ---
$slide6 = $p.Slides.Item(6)
$content6 = $slide6.Shapes.Item(2)
Set-TwoParagraphs $content6 "Vẽ liên tiếp nhiều đường tròn xếp cạnh nhau." "Mỗi đường tròn có bán kính khác nhau và giảm dần. "

# --- Slide 7: title text update ---
$slide7 = $p.Slides.Item(7)
$title7 = $slide7.Shapes.Item(1)
$title7.TextFrame.TextRange.Text = "6. Tô màu hình tròn (Cách 2)"

# --- Slide 7: content placeholder ---
$content7 = $slide7.Shapes.Item(2)
Set-TwoParagraphs $content7 "Vẽ nhiều đoạn thẳng xếp cạnh nhau." "Mỗi đoạn thẳng mới được bẻ góc bằng 1 độ."

# --- Update cached date field text (3/1/2023 -> 3/2/2023) on every layout
#     and on the slide master ---
$design = $p.Designs.Item(1)
$master = $design.SlideMaster

foreach ($shape in $master.Shapes) {
    if ($shape.Name -like "Date Placeholder*") {
        if ($shape.TextFrame.TextRange.Text -eq "3/1/2023") {
            $shape.TextFrame.TextRange.Text = "3/2/2023"
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    foreach ($shape in $layout.Shapes) {
        if ($shape.Name -like "Date Placeholder*") {
            if ($shape.TextFrame.TextRange.Text -eq "3/1/2023") {
                $shape.TextFrame.TextRange.Text = "3/2/2023"
            }
        }
    }
}
